$s0 = 'Std'
$s1 = 'hyperparameters = {
    "gamma": 0.9999, # 0.9977508883654606,
    "n_layers": 4,
    "h_size": 346,
    "dropout": 0.21887931565996233,
    "lr": 0.000643918092035234,
    "longevity_exponential": 1.04,
    "step_penalty_multiplier": 1.04,
    "ghost_reward": 0,
    "dot_extra_reward": 0,
    "energy_pill_extra_reward": 0,
    "optimizer": "SGD",
    "policy_file_name": "MsPacManPG_optimized.pt",
    "n_training_episodes": 100,
    "n_evaluation_episodes": 10,
    "max_t": 50000,
    "env_id": "ALE/MsPacman-ram-v5",
    "s_size": 128,
    "a_size": 5,
}'
$s2 = 'MsPacManPG_0.pt'
$s3 = 'Hyperparameters'
$s4 = 'Notes'
$s5 = 'Video'
$s6 = 'MsPacMan_replay_0.mp4'
$s7 = 'Best Avg'
$s8 = 'hyperparameters = {
    "gamma": 1, # 0.9977508883654606,
    "n_layers": 4,
    "h_size": 346,
    "dropout": 0.21887931565996233,
    "lr": 0.000643918092035234,
    "longevity_exponential": 1.04 ,  # 1.007765383540288,
    "step_penalty_multiplier": 1.04,
    "ghost_reward": 0,
    "dot_extra_reward": 0,
    "energy_pill_extra_reward": 0,
    "optimizer": "SGD",
    "policy_file_name": "MsPacManPG_optimized.pt",
    "n_training_episodes": 100,
    "n_evaluation_episodes": 10,
    "max_t": 50000,
    "env_id": "ALE/MsPacman-ram-v5",
    "s_size": 128,
    "a_size": 5,
}'
$s9 = 'Explores more, gets stuck in corner at the end'
$s10 = 'Explores a bit, gets stuck in corner'
$s11 = 'MsPacMan_replay_1.mp4'
$s12 = 'Last Avg'
$s13 = 'MsPacMan_replay_2.mp4'
$s14 = 'Hides in lower right corner. Scoring high from ghosts'
$s15 = 'MsPacManPG_2.pt'
$s16 = 'hyperparameters = {    
    ''gamma'': 0.9999773350962325, 
    ''n_layers'': 2, 
    ''h_size'': 180, 
    ''dropout'': 0.5030046989416794, 
    ''lr'': 0.0014691717217639733, 
    ''longevity_exponential'': 1.006105611179948, 
    ''step_penalty_multiplier'': 1.0193820217954723, 
    "ghost_reward": 0,
    ''dot_extra_reward'': 546, 
    ''energy_pill_extra_reward'': 964,
    
    "optimizer": "SGD",
    "policy_file_name": "MsPacManPG_optimized.pt",
    "n_training_episodes": 100,
    "n_evaluation_episodes": 10,
    "max_t": 50000,
    "env_id": "ALE/MsPacman-ram-v5",
    "s_size": 128,
    "a_size": 5,
}'
$s17 = 'Episodes of Training Done'
$s18 = 'Policy file name'
$s19 = 'MsPacManPG_optimized.pt'
$s20 = 'hyperparameters = {
    ''gamma'': 0.9999819193245816, 
    ''n_layers'': 1, 
    ''h_size'': 175, 
    ''dropout'': 0.44984866197635065, 
    ''lr'': 6.166629462708628e-05, 
    ''longevity_exponential'': 1.006491852944776, 
    ''step_penalty_multiplier'': 1.0386448544834312, 
    ''dot_extra_reward'': 13, 
    ''energy_pill_extra_reward'': 12,
    "ghost_reward": 0,
    "optimizer": "SGD",
    "policy_file_name": "MsPacManPG_optimized.pt",
    "n_training_episodes": 10000,
#     "n_evaluation_episodes": 10,
    "max_t": 50000,
    "env_id": "ALE/MsPacman-ram-v5",
    "s_size": 128,
    "a_size": 5,
}'
$s21 = '- MsPacMan_replay_3_best.mp4
- MsPacMan_replay_3_worst.mp4'
$s22 = 'hyperparameters_1 = {
    "h_size": 32,
    "n_training_episodes": 10000,
    "n_evaluation_episodes": 10,
    "max_t": 5000,
    "gamma": 0.99,
    "lr": 1e-4,
    "env_id": env_id,
    "state_space": 128,
    "action_space": 5,
}'
$s23 = 'MsPacMan_beforeOptimizing_replay.mp4'
$s24 = 'Final trained agent with improved policy network, rewards, and hyperparameters'
$s25 = 'Trained agent before making improvements to policy network, rewards, and hyperparameters'

$wb = $excel.ActiveWorkbook

# Remove the "Optuna Study" sheet entirely
$wsOptuna = $wb.Worksheets.Item("Optuna Study")
$wsOptuna.Delete()

$ws = $wb.Worksheets.Item("Policy")

# Wipe the existing table completely (content + formatting) so we can rebuild
# it cleanly at the new row positions without leaking stale styles/values.
$ws.Range("A1:I4").Clear()

# Row 1 - header (unchanged from before)
$ws.Range("A1").Value = $s18
$ws.Range("B1").Value = $s17
$ws.Range("C1").Value = $s7
$ws.Range("D1").Value = $s0
$ws.Range("E1").Value = $s12
$ws.Range("F1").Value = $s0
$ws.Range("G1").Value = $s3
$ws.Range("H1").Value = $s5
$ws.Range("I1").Value = $s4

# Row 2 - new "before optimizing" baseline run
$ws.Range("B2").Value = 10000
$ws.Range("G2").Value = $s22
$ws.Range("H2").Value = $s23
$ws.Range("I2").Value = $s25

# Row 3 - previously row 2
$ws.Range("A3").Value = $s2
$ws.Range("B3").Value = 100
$ws.Range("C3").Value = 943
$ws.Range("D3").Value = 847.18
$ws.Range("G3").Value = $s1
$ws.Range("H3").Value = $s6
$ws.Range("I3").Value = $s9

# Row 4 - previously row 3
$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 864
$ws.Range("D4").Value = 566.71
$ws.Range("G4").Value = $s8
$ws.Range("H4").Value = $s11
$ws.Range("I4").Value = $s10

# Row 5 - previously row 4
$ws.Range("A5").Value = $s15
$ws.Range("B5").Value = 100
$ws.Range("C5").Value = 712
$ws.Range("D5").Value = 447.83
$ws.Range("E5").Value = 435
$ws.Range("F5").Value = 206.94
$ws.Range("G5").Value = $s16
$ws.Range("H5").Value = $s13
$ws.Range("I5").Value = $s14

# Row 6 - new "final optimized" run
$ws.Range("A6").Value = $s19
$ws.Range("B6").Value = 10000
$ws.Range("G6").Value = $s20
$ws.Range("H6").Value = $s21
$ws.Range("I6").Value = $s24

# Selection ends up on G1 in the saved workbook
$ws.Range("G1").Select()
